$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 corresponds to file_name "metrics_sim_with_priors.json" (B3)
# Update metric values per corrected relevance markers.

$ws.Range("C3").Value = 1
$ws.Range("H3").Value = 0.9022711631108052
$ws.Range("I3").Value = 0.003491854032380856
$ws.Range("J3").Value = 0.8823529411764706
$ws.Range("K3").Value = 29.23529411764706

$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = 4
$ws.Range("T3").Value = 17
$ws.Range("U3").Value = 40

$ws.Range("V3").Value = 5795
$ws.Range("W3").Value = 5794
$ws.Range("X3").Value = 5791
$ws.Range("Y3").Value = 5778
$ws.Range("Z3").Value = 5755

$ws.Range("AF3").Value = 1
$ws.Range("AG3").Value = 0.999827
$ws.Range("AH3").Value = 0.99931
$ws.Range("AI3").Value = 0.997066
$ws.Range("AJ3").Value = 0.993097

$wb.Save()
